$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J12").Value = 150
$ws.Range("L12").Value = 150
$ws.Range("N12").Value = -490

$ws.Range("H19").Value = 890.3333
$ws.Range("I19").Value = 887.25
$ws.Range("J19").Value = 896.5
$ws.Range("K19").Value = 887.25
$ws.Range("L19").Value = 896.5
$ws.Range("M19").Value = -712.25
$ws.Range("N19").Value = -1246.5

$ws.Range("H41").Value = 359.42856
$ws.Range("I41").Value = 369.25
$ws.Range("J41").Value = 346.33334
$ws.Range("K41").Value = 369.25
$ws.Range("L41").Value = 346.33334
$ws.Range("M41").Value = 70.75
$ws.Range("N41").Value = -1226.33334

$ws.Range("H43").Value = 900
$ws.Range("I43").Value = 850
$ws.Range("J43").Value = 1000
$ws.Range("K43").Value = 850
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = -781
$ws.Range("N43").Value = -1138

$ws.Range("H62").Value = 8521.5
$ws.Range("I62").Value = 8521.5
$ws.Range("K62").Value = 8521.5
$ws.Range("M62").Value = -7897.5

$ws.Range("H65").Value = 8521.5
$ws.Range("I65").Value = 8521.5
$ws.Range("K65").Value = 42607.5
$ws.Range("M65").Value = -39487.5

$ws.Range("H76").Value = 3207
$ws.Range("I76").Value = 2716.6667
$ws.Range("K76").Value = 2716.6667
$ws.Range("M76").Value = -2401.6667

$ws.Range("H79").Value = 3207
$ws.Range("I79").Value = 2716.6667
$ws.Range("K79").Value = 2716.6667
$ws.Range("M79").Value = -1624.6667

$ws.Range("H82").Value = 8000
$ws.Range("I82").Value = 8000
$ws.Range("K82").Value = 24000
$ws.Range("M82").Value = -23594

$ws.Range("H85").Value = 8000
$ws.Range("I85").Value = 8000
$ws.Range("K85").Value = 24000
$ws.Range("M85").Value = -22596

$ws.Range("H92").Value = 709.12
$ws.Range("I92").Value = 816.3889
$ws.Range("J92").Value = 433.2857
$ws.Range("K92").Value = 816.3889
$ws.Range("L92").Value = 433.2857
$ws.Range("M92").Value = 431.6111
$ws.Range("N92").Value = -2929.2857

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H107").Value = 2412.4285
$ws.Range("I107").Value = 2412.4285
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2412.4285
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -492.4285
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = 754

$ws.Range("H115").Value = 1400
$ws.Range("I115").Value = 1400
$ws.Range("K115").Value = 4200
$ws.Range("M115").Value = -2633

$ws.Range("H116").Value = 4998
$ws.Range("J116").Value = 4997.5
$ws.Range("L116").Value = 4997.5
$ws.Range("N116").Value = -11881.5

$ws.Range("H137").Value = 1089.6666
$ws.Range("I137").Value = 866.625
$ws.Range("K137").Value = 2599.875
$ws.Range("M137").Value = -49.875

$ws.Range("H138").Value = 3940.1765
$ws.Range("J138").Value = 4043.6553
$ws.Range("L138").Value = 12130.9659
$ws.Range("N138").Value = -22410.9659

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1463.4615
$ws.Range("I2").Value = 1463.4615
$ws.Range("K2").Value = 1463.4615
$ws.Range("M2").Value = -1350.4615

$ws.Range("H32").Value = 18666.834
$ws.Range("I32").Value = 18000.176
$ws.Range("K32").Value = 18000.176
$ws.Range("M32").Value = -17713.176

$ws.Range("H88").Value = 2268.75
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 2450
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 2450
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -3262

$ws.Range("H91").Value = 2268.75
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 2450
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 2450
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -5258

$ws.Range("H97").Value = 2532.182
$ws.Range("I97").Value = 285.4
$ws.Range("K97").Value = 285.4
$ws.Range("M97").Value = 210.6

$ws.Range("H116").Value = 1463.4615
$ws.Range("I116").Value = 1463.4615
$ws.Range("K116").Value = 1463.4615
$ws.Range("M116").Value = 830.5385000000001

$ws.Range("H122").Value = 2597.3157
$ws.Range("I122").Value = 2412.7058
$ws.Range("K122").Value = 7238.117400000001
$ws.Range("M122").Value = -4788.117400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 669.1177
$ws.Range("I94").Value = 625.8
$ws.Range("K94").Value = 625.8
$ws.Range("M94").Value = -174.8

$ws.Range("H105").Value = 2565.182
$ws.Range("I105").Value = 2596.5557
$ws.Range("J105").Value = 2424
$ws.Range("K105").Value = 2596.5557
$ws.Range("L105").Value = 2424
$ws.Range("M105").Value = -849.5556999999999
$ws.Range("N105").Value = -5918

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2673.75
$ws.Range("I137").Value = 2673.75
$ws.Range("K137").Value = 8021.25
$ws.Range("M137").Value = -2921.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 35000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 35000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 35000
$ws.Range("N38").Value = -35946
$ws.Range("M38").ClearContents()

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H96").Value = 1749.75
$ws.Range("J96").Value = 1600
$ws.Range("L96").Value = 1600
$ws.Range("N96").Value = -4346

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()
